$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (6) to hold "Groupes CM"
$ws.Columns.Item(6).Insert()

# Header
$ws.Range("F1").Value = "Groupes CM"

# Fill value 1 for every data row (2..18)
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Re-apply the sort state so its stored ref grows from A2:I33 to A2:J33
# (Excel keeps remembering the last sort dialog's range over the whole table)
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B2:B33"))
$sort.SortFields.Add($ws.Range("A2:A33"))
$sort.SetRange($ws.Range("A2:J33"))
$sort.Apply()

# Update the selection to match the target state
$ws.Range("F18").Select()
